$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Update risk-register cell content (rows 8-17, plus header L7)
# -----------------------------------------------------------------
$ws.Range("L7").Value = "Proximity / Length of impact"

# --- Row 8 ---
$ws.Range("C8").Value = "Illness/Injury"
$ws.Range("D8").Value = "Developer lacks mobility and is unable to start or continue or finish the project."
$ws.Range("F8").Value = "M"
$ws.Range("I8").Value = "M"
$ws.Range("K8").Value = "Workload bulds-up and if developer does not recover before project deadline, project is incomplete."
$ws.Range("L8").Value = "Short term but variable depending on illness."
$ws.Range("M8").Value = "Avoid or limit human contact to bare necessities."

# --- Row 9 ---
$ws.Range("D9").Value = "Time allocated to each aspect of the project not clearly defined."
$ws.Range("F9").Value = "M"
$ws.Range("J9").Value = "H"
$ws.Range("K9").Value = "Overlap of individual feature timelines which creates a backlog and delay in project completion."
$ws.Range("L9").Value = "Long term."
$ws.Range("M9").Value = "Create a realistic deadline for each feature to be designed."

# --- Row 10 ---
$ws.Range("D10").Value = "Hardware performance is insufficient due to faulty or damaged equipment, i.e. laptop, computer,etc."
$ws.Range("E10").Value = "L"
$ws.Range("J10").Value = "H"
$ws.Range("K10").Value = "Fall behind deadlines due to unavailability of hardware to do work. Possible loss of data if hardware damaged."
$ws.Range("L10").Value = "Medium term issue and instance of occurrence unpredictable."
$ws.Range("M10").Value = "Contingency: continuous back-ups of project to avoid data loss and safeguard hardware to minimise risk of damage."

# --- Row 11 ---
$ws.Range("D11").Value = "Insufficient practice on programming language and mediocre knowledge of technologies to be used."
$ws.Range("F11").Value = "M"
$ws.Range("I11").Value = "M"
$ws.Range("K11").Value = "Low quality project developed, many of the requirements are not met."
$ws.Range("L11").Value = "Short term, initial stage of the project or for new implementations."
$ws.Range("M11").Value = "Research practices required for implementation of new features."

# --- Row 12 ---
$ws.Range("C12").Value = "Misunderstanding project scope"
$ws.Range("D12").Value = "Not reading and identifying the requirements correctly."
$ws.Range("E12").Value = "L"
$ws.Range("J12").Value = "H"
$ws.Range("K12").Value = "Objectives not met may result in redoing the project or failing."
$ws.Range("L12").Value = "Early stages, beginning of the project when objectives are stated."
$ws.Range("M12").Value = "Clarify requirements if not understood.`nProject completed to the given requirements."

# --- Row 13 ---
$ws.Range("B13").Value = "16.04.2020"
$ws.Range("D13").Value = "Errors in entity relationship diagram (ERD).`nRelationships between entities incorrectly defined."
$ws.Range("E13").Value = "L"
$ws.Range("I13").Value = "M"
$ws.Range("K13").Value = "Database and its tables created incorrectly."
$ws.Range("L13").Value = "Short term impact. Issue would be identified when inputting data to the databases tables."
$ws.Range("M13").Value = "Design an ERD which immitates the databases tables correctly and test it before advancing in the project."

# --- Row 14 ---
$ws.Range("B14").Value = "16.04.2020"
$ws.Range("D14").Value = "Low quality or no internet connection prevents accessing online cloud platforms such as GCP and nexus."
$ws.Range("E14").Value = "L"
$ws.Range("H14").Value = "L"
$ws.Range("K14").Value = "Unable to connect to database and view information in the tables to see if JAVA code works correctly. Unable to push files to nexus."
$ws.Range("L14").Value = "Short term issue,  but unpredictable."
$ws.Range("M14").Value = "Avoid: Reduce internet usage of other devices. Use mySQL workbench on device whenever internet is down."

# --- Row 15 ---
$ws.Range("B15").Value = "25.04.2020"
$ws.Range("C15").Value = "Git push "
$ws.Range("D15").Value = "Unable to upload files to the git repository."
$ws.Range("E15").Value = "L"
$ws.Range("H15").Value = "L"
$ws.Range("K15").Value = "Delays in uploading working features."
$ws.Range("L15").Value = "Short term."
$ws.Range("M15").Value = "Avoid: Only upload from a single device to avoid merge conflicts."

# --- Row 16 ---
$ws.Range("C16").Value = "Testing failures "
$ws.Range("D16").Value = "Low coverage of testing for the code developed."
$ws.Range("F16").Value = "M"
$ws.Range("J16").Value = "H"
$ws.Range("K16").Value = "Application not functioning as required so project requirements not met."
$ws.Range("L16").Value = "Medium term. "
$ws.Range("M16").Value = "Solve: Errors in code require refactoring to be able to pass the tests."

# --- Row 17 ---
$ws.Range("B17").Value = "29.04.2020"
$ws.Range("C17").Value = "Presentation questions"
$ws.Range("D17").Value = "Not well prepared to answer questions."
$ws.Range("E17").Value = "L"
$ws.Range("H17").Value = "L"
$ws.Range("K17").Value = "Understanding of project appears to be low and unable to explain the project to a high standard."
$ws.Range("L17").Value = "Short term, occuring at the end of the project."
$ws.Range("M17").Value = "Recap everything done to date to reaffirm understanding."


# -----------------------------------------------------------------
# 2. Give the new "Testing failures" row (16) its own date-formatted
#    cell for column B, matching the other register rows' look
#    (Arial Narrow 9pt, thin border, centered) but as a short date.
# -----------------------------------------------------------------
$refFont = $ws.Range("B8").Font
$b16 = $ws.Range("B16")
$b16.Value = "27.04.2020"
$b16.Font.Name = $refFont.Name
$b16.Font.Size = $refFont.Size
$b16.Font.Bold = $refFont.Bold
$b16.Font.Italic = $refFont.Italic
$b16.Font.Color = $refFont.Color
$b16.HorizontalAlignment = -4108
$b16.VerticalAlignment = -4108
$b16.Borders.LineStyle = 1
$b16.Borders.Weight = 2
$b16.NumberFormat = "mm-dd-yy"

# -----------------------------------------------------------------
# 3. Remove the now-unused trailing row (11th risk row, row 18) and
#    re-apply the shared "+1" counter formula over the smaller range
# -----------------------------------------------------------------
$ws.Rows.Item(18).Delete()
$ws.Range("A10:A17").Formula = "=A9+1"

# -----------------------------------------------------------------
# 4. Update the view: no frozen/scrolled top-left cell anymore, and
#    the active selection moves to M19 (just past the data)
# -----------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("M19").Select()

# -----------------------------------------------------------------
# 5. Touch row 2 so the sheet's recorded dimension keeps including
#    it (it carries only row-level formatting, no cell values)
# -----------------------------------------------------------------
$ws.Range("A2").Font.Bold = $ws.Range("A2").Font.Bold
